# Tripadvisor New Orleans shard: add a "State" column to hotel_info and
# reorder the sheets so review_info comes before hotel_info.

$wb = $excel.ActiveWorkbook

# --- 1. hotel_info: insert a new "State" column between Hotel_Name and City ---
$hotelWs = $wb.Worksheets.Item("hotel_info")

# Shift City/Zip/... one column to the right, opening up column C.
$hotelWs.Range("C:C").EntireColumn.Insert()

# Populate the new column: header + the single data row's value.
$hotelWs.Range("C1").Value = "State"
$hotelWs.Range("C2").Value = "Louisiana"

# --- 2. Reorder worksheets: review_info should precede hotel_info ---
$reviewWs = $wb.Worksheets.Item("review_info")
$reviewWs.Move($hotelWs)
